$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - only B changes
$ws.Range("B24").Value = 88637

# Row 25
$ws.Range("A25").Value = 112013691
$ws.Range("B25").Value = 88637
$ws.Range("E25").Value = 1962
$ws.Range("F25").Value = "Vaddporing"
$ws.Range("G25").Value = "Anomoporia kamtschatica"
$ws.Range("H25").Value = "(Parmasto) Bondartseva"
$ws.Range("Q25").Value = 610134
$ws.Range("R25").Value = 7121461
$ws.Range("Z25").Value = "19:29"
$ws.Range("AB25").Value = "19:29"

# Row 26
$ws.Range("A26").Value = 112013696
$ws.Range("B26").Value = 87109
$ws.Range("E26").Value = 4962
$ws.Range("F26").Value = "Mjölsvärting"
$ws.Range("G26").Value = "Lyophyllum semitale"
$ws.Range("H26").Value = "(Fr. : Fr.) Kühner"
$ws.Range("Q26").Value = 610070
$ws.Range("R26").Value = 7121402
$ws.Range("Z26").Value = "19:40"
$ws.Range("AB26").Value = "19:40"

# Row 27 - only B changes
$ws.Range("B27").Value = 89571

# Row 28
$ws.Range("A28").Value = 112013698
$ws.Range("B28").Value = 77650
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 610094
$ws.Range("R28").Value = 7121456
$ws.Range("Z28").Value = "19:49"
$ws.Range("AB28").Value = "19:49"

# Row 29
$ws.Range("A29").Value = 112013699
$ws.Range("B29").Value = 77650
$ws.Range("Q29").Value = 610068
$ws.Range("R29").Value = 7121408
$ws.Range("Z29").Value = "19:40"
$ws.Range("AB29").Value = "19:40"

# Row 30
$ws.Range("A30").Value = 112013704
$ws.Range("B30").Value = 81385
$ws.Range("E30").Value = 1312
$ws.Range("F30").Value = "Gammelgransskål"
$ws.Range("G30").Value = "Pseudographis pinicola"
$ws.Range("H30").Value = "(Nyl.) Rehm"
$ws.Range("Q30").Value = 610094
$ws.Range("R30").Value = 7121455
$ws.Range("Z30").Value = "19:49"
$ws.Range("AB30").Value = "19:49"

# Row 31
$ws.Range("A31").Value = 112013690
$ws.Range("B31").Value = 88637
$ws.Range("E31").Value = 1962
$ws.Range("F31").Value = "Vaddporing"
$ws.Range("G31").Value = "Anomoporia kamtschatica"
$ws.Range("H31").Value = "(Parmasto) Bondartseva"
$ws.Range("Q31").Value = 610052
$ws.Range("R31").Value = 7121425
$ws.Range("Z31").Value = "19:43"
$ws.Range("AB31").Value = "19:43"

# Row 32
$ws.Range("A32").Value = 112013700
$ws.Range("B32").Value = 77650
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = "Garnlav"
$ws.Range("G32").Value = "Alectoria sarmentosa"
$ws.Range("H32").Value = "(Ach.) Ach."
$ws.Range("Q32").Value = 610102
$ws.Range("R32").Value = 7121416
$ws.Range("Z32").Value = "19:35"
$ws.Range("AB32").Value = "19:35"

# Row 33
$ws.Range("A33").Value = 112013703
$ws.Range("B33").Value = 77650
$ws.Range("Q33").Value = 610144
$ws.Range("R33").Value = 7121461
$ws.Range("Z33").Value = "19:28"
$ws.Range("AB33").Value = "19:28"
